$d = $word.ActiveDocument

$find = "Biblio doit contenir les 4 références : Mais aussi voir pour faire la biblio d’office en entier à partir du fichier bibtex (bio CV)."
$replace = "Biblio doit contenir les 4 références : Mais aussi voir pour faire la biblio d’office en entier à partir du fichier bibtex (bio CV)."

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
